# Insert a new data row at row 14 (pushes existing rows 14-62 down to 15-63)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with its data
$ws.Cells.Item(14, 1).Value = 10
$ws.Cells.Item(14, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(14, 3).Value = "La Araucanía"
$ws.Cells.Item(14, 4).Value = 44701
$ws.Cells.Item(14, 5).Value = 9
$ws.Cells.Item(14, 6).Value = 300000001
$ws.Cells.Item(14, 7).Value = "Rabanito"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 40
$ws.Cells.Item(14, 11).Value = 7000
$ws.Cells.Item(14, 12).Value = 7000
$ws.Cells.Item(14, 13).Value = 7000
$ws.Cells.Item(14, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(14, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(14, 16).Value = 583
$ws.Cells.Item(14, 17).Value = 12
$ws.Cells.Item(14, 18).Value = "Hortaliza"
